$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking strings
# (e.g. "1.000", "27.813.65") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row -> Price(D), Volume1h(E) updates. Only cells that actually change are set.
$updates = @(
    @{ Row = 2;  D = "27.813.65";      E = "  -2.91%  " },
    @{ Row = 3;  D = "1.791.55";       E = "  -0.75%  " },
    @{ Row = 4;                       E = "  -0.21%  " },
    @{ Row = 5;  D = "316.21";         E = "  -0.36%  " },
    @{ Row = 6;  D = "0.9999";         E = "  -0.28%  " },
    @{ Row = 7;  D = "0.5324";         E = "  -0.04%  " },
    @{ Row = 8;  D = "0.3846";         E = "  +1.56%  " },
    @{ Row = 9;                       E = "  -1.01%  " },
    @{ Row = 10; D = "41.43";          E = "  -2.61%  " },
    @{ Row = 11;                       E = "  -2.70%  " },
    @{ Row = 12;                       E = "  -0.18%  " },
    @{ Row = 13; D = "6.179";          E = "  +0.14%  " },
    @{ Row = 14; D = "7.466";          E = "  +1.50%  " },
    @{ Row = 15; D = "20.32";          E = "  -1.90%  " },
    @{ Row = 16; D = "1.791.81";       E = "  -0.45%  " },
    @{ Row = 17; D = "88.22";          E = "  -2.48%  " },
    @{ Row = 18; D = "0.00001058";     E = "  -0.94%  " },
    @{ Row = 19; D = "0.06522";        E = "  +1.17%  " },
    @{ Row = 20; D = "1.000";          E = "  -0.19%  " },
    @{ Row = 21; D = "17.20";          E = "  -0.22%  " },
    @{ Row = 22; D = "5.954";          E = "  +0.48%  " },
    @{ Row = 23; D = "27.870.27";      E = "  -2.73%  " },
    @{ Row = 24; D = "11.15";          E = "  +0.76%  " },
    @{ Row = 25; D = "2.094";          E = "  -0.18%  " },
    @{ Row = 26; D = "157.08";         E = "  -1.92%  " },
    @{ Row = 27; D = "20.14";          E = "  -1.59%  " },
    @{ Row = 28; D = "1.997.57";       E = "  -0.59%  " },
    @{ Row = 29;                       E = "  -2.51%  " },
    @{ Row = 30; D = "121.38";         E = "  -1.21%  " },
    @{ Row = 31; D = "0.1092";         E = "  +4.23%  " },
    @{ Row = 32; D = "1.100";          E = "  -0.38%  " },
    @{ Row = 33; D = "3.652";          E = "  -0.89%  " },
    @{ Row = 34;                       E = "  -2.78%  " },
    @{ Row = 35; D = "0.06906";        E = "  +7.34%  " },
    @{ Row = 36; D = "0.2197";         E = "  -2.63%  " },
    @{ Row = 37; D = "0.02261";        E = "  -2.24%  " },
    @{ Row = 38; D = "5.042";          E = "  +0.02%  " },
    @{ Row = 39; D = "11.43";          E = "  +1.13%  " },
    @{ Row = 40; D = "8.377";          E = "  -5.00%  " },
    @{ Row = 41; D = "0.6095";         E = "  -2.32%  " },
    @{ Row = 42; D = "1.169";          E = "  -5.22%  " },
    @{ Row = 43;                       E = "  +0.49%  " },
    @{ Row = 44; D = "13.22";          E = "  -0.93%  " },
    @{ Row = 45; D = "3.676";          E = "  -0.59%  " },
    @{ Row = 46; D = "0.5684";         E = "  -3.14%  " },
    @{ Row = 47; D = "125.00";         E = "  -0.83%  " },
    @{ Row = 48; D = "1.908";          E = "  -2.16%  " },
    @{ Row = 49; D = "1.169" },
    @{ Row = 50; D = "0.06795";        E = "  -1.39%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("D")) {
        $ws.Cells.Item($r, 4).Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}

# Row 51 is a full row replacement: Aave -> BabyDogeCoin
$ws.Cells.Item(51, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(51, 4).Value = "0.00000000291"
$ws.Cells.Item(51, 5).Value = "  +36.54%  "
